$d = $word.ActiveDocument

# Locate the word "query" in the first paragraph (right after the
# spellStart proofErr mark) so we can insert the warning text just
# before it.
$found = $d.Content.Find.Execute("query", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $d.Content.Start

# NOTE: Find.Execute above mutates $d.Content itself only when called on
# that exact range object; recompute explicitly with a fresh search range.
$searchRange = $d.Content
$searchRange.Find.Execute("query", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $searchRange.Start

# 1) Leading "    " (four spaces), plain formatting.
$r = $d.Range($pos, $pos)
$r.InsertBefore("    ")
$pos = $pos + 4

# 2) "<---" marker, colored/sized/highlighted.
$markerText = "<---"
$r = $d.Range($pos, $pos)
$r.InsertBefore($markerText)
$markerRange = $d.Range($pos, $pos + $markerText.Length)
$markerRange.Font.Color = 42495
$markerRange.Font.Size = 16
$markerRange.Font.HighlightColorIndex = 16
$pos = $pos + $markerText.Length

# 3) Version mismatch message, same formatting as the marker.
$msgText = "M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0"
$r = $d.Range($pos, $pos)
$r.InsertBefore($msgText)
$msgRange = $d.Range($pos, $pos + $msgText.Length)
$msgRange.Font.Color = 42495
$msgRange.Font.Size = 16
$msgRange.Font.HighlightColorIndex = 16
$pos = $pos + $msgText.Length

# 4) Trailing "    " (four spaces), plain formatting.
$r = $d.Range($pos, $pos)
$r.InsertBefore("    ")
$pos = $pos + 4
